$wb = $excel.ActiveWorkbook

# Sheet "展览" updates
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F6").Value = 851
$ws1.Range("F8").Value = 735
$ws1.Range("F9").Value = 13157
$ws1.Range("F10").Value = 13030
$ws1.Range("F15").Value = 66
$ws1.Range("F16").Value = 632
$ws1.Range("F17").Value = 2059

# Sheet "全部类型" updates
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F7").Value = 851
$ws4.Range("F10").Value = 735
$ws4.Range("F11").Value = 13157
$ws4.Range("F12").Value = 13030
$ws4.Range("F17").Value = 66
$ws4.Range("F18").Value = 632
$ws4.Range("F21").Value = 2059
